# =============================================================================
# initialize-data.xlsx edit script
#
#  1. Rename Sheet3 -> TransactionStatuses and populate it with a header +
#     2 data rows (Completed / Cancelled) mirroring the Permissions sheet's
#     "delete/insert" SQL-script generator pattern.
#  2. Permissions sheet:
#       - fix a handful of Serbian-Latin translations (NameLatin column)
#       - drop the stray "Insert users" permission (rows simply get
#         re-purposed/overwritten below - the row that used to hold it now
#         holds "Read permissions")
#       - append 9 new permission rows (Permissions/Tiers/Transactions CRUD)
#       - append a handful of blank formatted rows at the bottom
# =============================================================================

$wb = $excel.ActiveWorkbook

$home        = $wb.Worksheets.Item(1)
$permissions = $wb.Worksheets.Item(2)
$txStatuses  = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------------
# 1) Rename Sheet3 -> TransactionStatuses
# ---------------------------------------------------------------------------
$txStatuses.Name = "TransactionStatuses"

# ---------------------------------------------------------------------------
# 2) Permissions sheet - helper to build the shared "insert into ..." formula
#    text for a given data row (same shape as the existing H3 formula, just
#    re-anchored row-wise)
# ---------------------------------------------------------------------------
function New-PermissionsFormula([int]$r) {
    return '=CONCATENATE("insert into ",$A$1,"(",$B$2,", ",$C$2,", ",$D$2,", ",$E$2,", ",$F$2,", ",$G$2,") values(N''",B' + $r + ',"'', N''",C' + $r + ',"'', ",IF(TRIM(D' + $r + ')<>"","N''"&D' + $r + '&"''","null"),", ",IF(TRIM(E' + $r + ')<>"","N''"&E' + $r + '&"''","null")," , getdate(), N''",G' + $r + ',"'');")'
}

# Full B (Name) / C (NameLatin) / G (Code) content for data rows 3-17.
# Rows 3-8 keep their original "Name" (English) text but get corrected
# "NameLatin" translations; row 9 is completely re-purposed (it used to be
# "Insert users", now holds what used to live in row 10 - "Delete users");
# rows 10-17 are brand new permission rows.
$permRows = @(
    @{ Row = 3;  B = "Read roles";                 C = "Pregled uloga korisnika";            G = "ReadRole" }
    @{ Row = 4;  B = "Edit roles";                 C = "Promena uloga korisnika";            G = "EditRole" }
    @{ Row = 5;  B = "Insert roles";                C = "Dodavanje uloga korisnika";          G = "InsertRole" }
    @{ Row = 6;  B = "Delete roles";                C = "Brisanje uloga korisnika";           G = "DeleteRole" }
    @{ Row = 7;  B = "Read users";                  C = "Pregled profila korisnika";          G = "ReadUserExtended" }
    @{ Row = 8;  B = "Edit users";                  C = "Promena profila korisnika";          G = "EditUserExtended" }
    @{ Row = 9;  B = "Delete users";                C = "Brisanje profila korisnika";         G = "DeleteUserExtended" }
    @{ Row = 10; B = "Read permissions";            C = "Pregled permisija uloga";            G = "ReadPermission" }
    @{ Row = 11; B = "Read tiers";                  C = "Pregled nivoa odanosti";             G = "ReadTier" }
    @{ Row = 12; B = "Edit tiers";                  C = "Promena nivoa odanosti";             G = "EditTier" }
    @{ Row = 13; B = "Insert tiers";                C = "Dodavanje nivoa odanosti";           G = "InsertTier" }
    @{ Row = 14; B = "Delete tiers";                C = "Brisanje nivoa odanosti";            G = "DeleteTier" }
    @{ Row = 15; B = "Read transactions";           C = "Pregled transakcija";                G = "ReadTransaction" }
    @{ Row = 16; B = "Read transaction products";   C = "Pregled proizvoda iz transakcije";   G = "ReadTransactionProduct" }
    @{ Row = 17; B = "Read transaction statuses";   C = "Pregled statusa transakcije";        G = $null }
)

foreach ($pr in $permRows) {
    $r = $pr.Row
    $permissions.Range("A$r").Value2 = ($r - 2)
    $permissions.Range("B$r").Value2 = $pr.B
    $permissions.Range("C$r").Value2 = $pr.C
    if ($pr.G -ne $null) {
        $permissions.Range("G$r").Value2 = $pr.G
        $permissions.Range("H$r").Formula = (New-PermissionsFormula $r)
    }
}

# Rows 18-23: trailing blank (but formatted) rows
for ($r = 18; $r -le 23; $r++) {
    $permissions.Range("A$r`:C$r").Value2 = ""
}

# ---------------------------------------------------------------------------
# 4) Column widths / view state on the Permissions sheet
# ---------------------------------------------------------------------------
$permissions.Columns.Item(2).ColumnWidth = 26.166666666666668
$permissions.Columns.Item(3).ColumnWidth = 30.833333333333332
$permissions.Range("C21").Select()

# ---------------------------------------------------------------------------
# 5) TransactionStatuses sheet content
# ---------------------------------------------------------------------------
$txStatuses.Range("A1").Value2 = "TransactionStatuses"

$txStatuses.Range("A2").Value2 = "Id"
$txStatuses.Range("B2").Value2 = "Name"
$txStatuses.Range("C2").Value2 = "NameLatin"
$txStatuses.Range("D2").Value2 = "Code"
$txStatuses.Range("E2").Formula = '=CONCATENATE("delete from ",$A$1,"; dbcc checkident (",$A$1,", reseed, 0);")'

$txStatuses.Range("A3").Value2 = 1
$txStatuses.Range("B3").Value2 = "Completed"
$txStatuses.Range("C3").Value2 = "Kompletirana"
$txStatuses.Range("D3").Value2 = "Completed"
$txStatuses.Range("E3").Formula = '=CONCATENATE("insert into ",$A$1,"(",$B$2,", ",$C$2,", ",$D$2,") values(N''",B3,"'', N''",C3,"'', N''",D3,"'');")'

$txStatuses.Range("A4").Value2 = 2
$txStatuses.Range("B4").Value2 = "Cancelled"
$txStatuses.Range("C4").Value2 = "Otkazana"
$txStatuses.Range("D4").Value2 = "Cancelled"
$txStatuses.Range("E4").Formula = '=CONCATENATE("insert into ",$A$1,"(",$B$2,", ",$C$2,", ",$D$2,") values(N''",B4,"'', N''",C4,"'', N''",D4,"'');")'

# Column widths / view state on the TransactionStatuses sheet
$txStatuses.Columns.Item(2).ColumnWidth = 22.666666666666668
$txStatuses.Columns.Item(3).ColumnWidth = 24.833333333333332
$txStatuses.Columns.Item(4).ColumnWidth = 12.666666666666666
$txStatuses.Range("D3:D4").Select()

Write-Output "edit complete"
